$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 100.8744713333333
$ws.Range("H2").Value = 302.623414
$ws.Range("I2").Value = 0.1452075237922473
$ws.Range("J2").Value = 0.1452075237922473
$ws.Range("M2").Value = 2.759544333333333
$ws.Range("N2").Value = 8.278632999999999
$ws.Range("O2").Value = 0.2574067337278401
$ws.Range("P2").Value = 0.2574067337278401
$ws.Range("Q2").Value = 278.3675757458957
$ws.Range("R2").Value = 2505.308181713062
$ws.Range("S2").Value = 0.03737739441207001
$ws.Range("T2").Value = 0.03737739441207001

$ws.Range("G3").Value = 100.8744713333333
$ws.Range("H3").Value = 302.623414
$ws.Range("I3").Value = 0.1452075237922473
$ws.Range("J3").Value = 0.1452075237922473
$ws.Range("O3").Value = 0.6758254232987829
$ws.Range("P3").Value = 0.6758254232987829
$ws.Range("Q3").Value = 730.8584433149928
$ws.Range("R3").Value = 6577.725989834935
$ws.Range("S3").Value = 0.09813493623306364
$ws.Range("T3").Value = 0.09813493623306364

$ws.Range("G4").Value = 100.8744713333333
$ws.Range("H4").Value = 302.623414
$ws.Range("I4").Value = 0.1452075237922473
$ws.Range("J4").Value = 0.1452075237922473
$ws.Range("M4").Value = 0.5200313333333334
$ws.Range("N4").Value = 1.560094
$ws.Range("O4").Value = 0.0485078515798926
$ws.Range("P4").Value = 0.0485078515798926
$ws.Range("Q4").Value = 52.45788582676846
$ws.Range("R4").Value = 472.1209724409161
$ws.Range("S4").Value = 0.007043705012398057
$ws.Range("T4").Value = 0.007043705012398056

$ws.Range("G5").Value = 100.8744713333333
$ws.Range("H5").Value = 302.623414
$ws.Range("I5").Value = 0.1452075237922473
$ws.Range("J5").Value = 0.1452075237922473
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1957573333333333
$ws.Range("N5").Value = 0.587272
$ws.Range("O5").Value = 0.01825999139348442
$ws.Range("P5").Value = 0.01825999139348442
$ws.Range("Q5").Value = 19.74691750962311
$ws.Range("R5").Value = 177.722257586608
$ws.Range("S5").Value = 0.002651488134715621
$ws.Range("T5").Value = 0.002651488134715621

$ws.Range("I6").Value = 0.7769829249672668
$ws.Range("J6").Value = 0.776982924967267
$ws.Range("M6").Value = 2.759544333333333
$ws.Range("N6").Value = 8.278632999999999
$ws.Range("O6").Value = 0.2574067337278401
$ws.Range("P6").Value = 0.2574067337278401
$ws.Range("Q6").Value = 1489.501697780766
$ws.Range("R6").Value = 13405.51528002689
$ws.Range("S6").Value = 0.2000006368781276
$ws.Range("T6").Value = 0.2000006368781276

$ws.Range("I7").Value = 0.7769829249672668
$ws.Range("J7").Value = 0.776982924967267
$ws.Range("O7").Value = 0.6758254232987829
$ws.Range("P7").Value = 0.6758254232987829
$ws.Range("S7").Value = 0.5251048141619296
$ws.Range("T7").Value = 0.5251048141619297

$ws.Range("I8").Value = 0.7769829249672668
$ws.Range("J8").Value = 0.776982924967267
$ws.Range("M8").Value = 0.5200313333333334
$ws.Range("N8").Value = 1.560094
$ws.Range("O8").Value = 0.0485078515798926
$ws.Range("P8").Value = 0.0485078515798926
$ws.Range("Q8").Value = 280.6940060874285
$ws.Range("R8").Value = 2526.246054786856
$ws.Range("S8").Value = 0.03768977240442301
$ws.Range("T8").Value = 0.03768977240442301

$ws.Range("I9").Value = 0.7769829249672668
$ws.Range("J9").Value = 0.776982924967267
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1957573333333333
$ws.Range("N9").Value = 0.587272
$ws.Range("O9").Value = 0.01825999139348442
$ws.Range("P9").Value = 0.01825999139348442
$ws.Range("Q9").Value = 105.6626910577031
$ws.Range("R9").Value = 950.964219519328
$ws.Range("S9").Value = 0.01418770152278664
$ws.Range("T9").Value = 0.01418770152278665

$ws.Range("G10").Value = 53.798087
$ws.Range("H10").Value = 161.394261
$ws.Range("I10").Value = 0.07744166482137986
$ws.Range("J10").Value = 0.07744166482137986
$ws.Range("M10").Value = 2.759544333333333
$ws.Range("N10").Value = 8.278632999999999
$ws.Range("O10").Value = 0.2574067337278401
$ws.Range("P10").Value = 0.2574067337278401
$ws.Range("Q10").Value = 148.4582061250237
$ws.Range("R10").Value = 1336.123855125213
$ws.Range("S10").Value = 0.01993400599611756
$ws.Range("T10").Value = 0.01993400599611756

$ws.Range("G11").Value = 53.798087
$ws.Range("H11").Value = 161.394261
$ws.Range("I11").Value = 0.07744166482137986
$ws.Range("J11").Value = 0.07744166482137986
$ws.Range("O11").Value = 0.6758254232987829
$ws.Range("P11").Value = 0.6758254232987829
$ws.Range("Q11").Value = 389.779352480749
$ws.Range("R11").Value = 3508.014172326741
$ws.Range("S11").Value = 0.05233704590887151
$ws.Range("T11").Value = 0.05233704590887151

$ws.Range("G12").Value = 53.798087
$ws.Range("H12").Value = 161.394261
$ws.Range("I12").Value = 0.07744166482137986
$ws.Range("J12").Value = 0.07744166482137986
$ws.Range("M12").Value = 0.5200313333333334
$ws.Range("N12").Value = 1.560094
$ws.Range("O12").Value = 0.0485078515798926
$ws.Range("P12").Value = 0.0485078515798926
$ws.Range("Q12").Value = 27.97669091339267
$ws.Range("R12").Value = 251.790218220534
$ws.Range("S12").Value = 0.003756528783255285
$ws.Range("T12").Value = 0.003756528783255284

$ws.Range("G13").Value = 53.798087
$ws.Range("H13").Value = 161.394261
$ws.Range("I13").Value = 0.07744166482137986
$ws.Range("J13").Value = 0.07744166482137986
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1957573333333333
$ws.Range("N13").Value = 0.587272
$ws.Range("O13").Value = 0.01825999139348442
$ws.Range("P13").Value = 0.01825999139348442
$ws.Range("Q13").Value = 10.53137004955467
$ws.Range("R13").Value = 94.782330445992
$ws.Range("S13").Value = 0.001414084133135502
$ws.Range("T13").Value = 0.001414084133135502

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.2555676666666666
$ws.Range("H14").Value = 0.7667029999999999
$ws.Range("I14").Value = 0.0003678864191059829
$ws.Range("J14").Value = 0.000367886419105983
$ws.Range("M14").Value = 2.759544333333333
$ws.Range("N14").Value = 8.278632999999999
$ws.Range("O14").Value = 0.2574067337278401
$ws.Range("P14").Value = 0.2574067337278401
$ws.Range("Q14").Value = 0.7052503063332221
$ws.Range("R14").Value = 6.347252756998999
$ws.Range("S14").Value = [double]"9.469644152490231E-05"
$ws.Range("T14").Value = [double]"9.469644152490234E-05"

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.2555676666666666
$ws.Range("H15").Value = 0.7667029999999999
$ws.Range("I15").Value = 0.0003678864191059829
$ws.Range("J15").Value = 0.000367886419105983
$ws.Range("O15").Value = 0.6758254232987829
$ws.Range("P15").Value = 0.6758254232987829
$ws.Range("Q15").Value = 1.851645758860333
$ws.Range("R15").Value = 16.664811829743
$ws.Range("S15").Value = 0.0002486269949181743
$ws.Range("T15").Value = 0.0002486269949181744

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.2555676666666666
$ws.Range("H16").Value = 0.7667029999999999
$ws.Range("I16").Value = 0.0003678864191059829
$ws.Range("J16").Value = 0.000367886419105983
$ws.Range("M16").Value = 0.5200313333333334
$ws.Range("N16").Value = 1.560094
$ws.Range("O16").Value = 0.0485078515798926
$ws.Range("P16").Value = 0.0485078515798926
$ws.Range("Q16").Value = 0.1329031944535556
$ws.Range("R16").Value = 1.196128750082
$ws.Range("S16").Value = [double]"1.784537981625118E-05"
$ws.Range("T16").Value = [double]"1.784537981625119E-05"

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.2555676666666666
$ws.Range("H17").Value = 0.7667029999999999
$ws.Range("I17").Value = 0.0003678864191059829
$ws.Range("J17").Value = 0.000367886419105983
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1957573333333333
$ws.Range("N17").Value = 0.587272
$ws.Range("O17").Value = 0.01825999139348442
$ws.Range("P17").Value = 0.01825999139348442
$ws.Range("Q17").Value = 0.05002924491288888
$ws.Range("R17").Value = 0.450263204216
$ws.Range("S17").Value = [double]"6.717602846655051E-06"
$ws.Range("T17").Value = [double]"6.717602846655052E-06"
